$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K ("Hyperscaler"), shifting existing columns K:V to L:W.
# (format is copied from the left neighbour, which is what produces the s="9"/s="8"/s="1"
# styles seen on the new column in rows 1-5, matching a native Excel column insert)
$ws.Columns("K").Insert()

# Header text for the newly inserted column, and its (best-fit) width.
$ws.Range("K3").Value = "Hyperscaler"
$ws.Columns("K").ColumnWidth = 10.5

# The merged title bar (row1) / sub-header (row2) grew from A1:R1 / B2:P2 to A1:S1 / B2:Q2
# automatically because of the column insert; re-assert the selection to match.
$ws.Range("A1:S1").Select() | Out-Null

# Update the two values that changed in the (now shifted) Storage Node columns.
$ws.Range("U4").Value = "SG5812"
$ws.Range("U5").Value = "SG5812"
$ws.Range("V4").Value = "4TB HDDs Non-FDE"
$ws.Range("V5").Value = "4TB HDDs Non-FDE"
